$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L ("nemad" / ticker symbol) incorrectly held the company name
# ("نفت سپاهان"); update it to the correct ticker symbol ("شسپا") for
# every data row (2 through 45).
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 12).Value = "شسپا"
}

# Set column L width (Excel auto "best fit" after the text change) and
# select L2 as the active cell, matching the saved view state.
$ws.Columns.Item(12).ColumnWidth = 10
$ws.Range("L2").Select() | Out-Null
